$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 106 (shifts existing rows 106-117 down to 108-119)
$ws.Rows.Item(106).Resize(2).Insert()

# Expand the Excel Table ("Tableau2") to cover the two new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E119"))

# New row 106: MARTIN / Camille
$ws.Range("A106").Value = "MARTIN "
$ws.Range("B106").Value = "Camille"
$ws.Range("C106").Value = "N/A"
$ws.Range("D106").Value = "N/A"
$ws.Range("E106").Value = 1

# New row 107: TANCHOUX / Helene
$ws.Range("A107").Value = "TANCHOUX"
$ws.Range("B107").Value = "Helene"
$ws.Range("C107").Value = "N/A"
$ws.Range("D107").Value = "N/A"
$ws.Range("E107").Value = 1

# Match the "N/A" cell shading used elsewhere in the sheet (e.g. C5:D6) for the new rows
$ws.Range("C5:D6").Copy()
$ws.Range("C106:D107").PasteSpecial(-4122)

# Restore the active selection to reflect where the user ended up editing
$ws.Range("E108").Select()
